# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 606
$wsExhibit.Range("F7").Value = 2587
$wsExhibit.Range("F9").Value = 7035
$wsExhibit.Range("F11").Value = 442
$wsExhibit.Range("F13").Value = 93

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 606
$wsAll.Range("F9").Value = 2587
$wsAll.Range("F11").Value = 7035
$wsAll.Range("F13").Value = 442
$wsAll.Range("F17").Value = 93
